# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the 6c124c18-... file after a new handback round was processed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     6c124c18-...md row (row 2) ---
$overview = $wb.Sheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-16 18:42:21"

# --- zh-cn sheet: refresh handoff / handback datetimes for the
#     6c124c18-...md row (row 2) ---
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-16 18:42:17"
$zhcn.Range("K2").Value = "2016-08-16 18:42:33"

# --- de-de sheet: refresh handoff / handback datetimes for the
#     6c124c18-...md row (row 2) ---
$dede = $wb.Sheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-16 18:42:21"
$dede.Range("K2").Value = "2016-08-16 18:42:40"
